$d = $word.ActiveDocument

# Locate the "BILAGA 1 - Fridlysta arter" title paragraph; the new knärot
# section is appended directly after it (and before the trailing sectPr).
$anchorIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd() -eq "BILAGA 1 - Fridlysta arter") {
        $anchorIdx = $i
        break
    }
}
$anchor = $d.Paragraphs.Item($anchorIdx)
$rng = $anchor.Range
$rng.Collapse(0)

# Pre-create all 13 empty paragraphs first. Doing the InsertParagraphAfter
# calls up front (rather than interleaving with text insertion) sidesteps a runtime
# quirk where typing into a paragraph that is *currently* the last one in the document
# body mis-anchors the insertion point by one character.
for ($i = 0; $i -lt 13; $i++) {
    $rng.InsertParagraphAfter()
    $rng.Collapse(0)
    $p = $d.Paragraphs.Item($d.Paragraphs.Count)
    $rng = $p.Range
    $rng.Collapse(0)
}

# Paragraph 0: style=Heading1
$p = $d.Paragraphs.Item($anchorIdx + 1)
$p.Style = $d.Styles.Item('Heading 1')
$pos = $p.Range.Start
$r = $d.Range($pos, $pos)
$r.InsertAfter('Knärot – ekologi samt krav på livsmiljön')
$pos = $pos + 40

# Paragraph 1: style=None
$p = $d.Paragraphs.Item($anchorIdx + 2)
$p.Style = $d.Styles.Item("Normal")
$pos = $p.Range.Start
$r = $d.Range($pos, $pos)
$r.InsertAfter('Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).')
$pos = $pos + 684

# Paragraph 2: style=None
$p = $d.Paragraphs.Item($anchorIdx + 3)
$p.Style = $d.Styles.Item("Normal")
$pos = $p.Range.Start
$r = $d.Range($pos, $pos)
$r.InsertAfter('Samuel Johnsons doktorsavhandling ')
$pos = $pos + 34
$r = $d.Range($pos, $pos)
$r.InsertAfter('“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“')
$d.Range($pos, $pos + 82).Font.Italic = $true
$pos = $pos + 82
$r = $d.Range($pos, $pos)
$r.InsertAfter(' (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: ')
$pos = $pos + 162
$r = $d.Range($pos, $pos)
$r.InsertAfter('“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” ')
$d.Range($pos, $pos + 205).Font.Italic = $true
$pos = $pos + 205
$r = $d.Range($pos, $pos)
$r.InsertAfter('Vidare ')
$pos = $pos + 7
$r = $d.Range($pos, $pos)
$r.InsertAfter('“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”')
$d.Range($pos, $pos + 118).Font.Italic = $true
$pos = $pos + 118

# Paragraph 3: style=None
$p = $d.Paragraphs.Item($anchorIdx + 4)
$p.Style = $d.Styles.Item("Normal")
$pos = $p.Range.Start
$r = $d.Range($pos, $pos)
$r.InsertAfter('Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: ')
$pos = $pos + 205
$r = $d.Range($pos, $pos)
$r.InsertAfter('“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”')
$d.Range($pos, $pos + 865).Font.Italic = $true
$pos = $pos + 865

# Paragraph 4: style=None
$p = $d.Paragraphs.Item($anchorIdx + 5)
$p.Style = $d.Styles.Item("Normal")
$pos = $p.Range.Start
$r = $d.Range($pos, $pos)
$r.InsertAfter('En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).')
$pos = $pos + 1337

# Paragraph 5: style=None
$p = $d.Paragraphs.Item($anchorIdx + 6)
$p.Style = $d.Styles.Item("Normal")
$pos = $p.Range.Start
$r = $d.Range($pos, $pos)
$r.InsertAfter('Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).')
$pos = $pos + 868

# Paragraph 6: style=Heading2
$p = $d.Paragraphs.Item($anchorIdx + 7)
$p.Style = $d.Styles.Item('Heading 2')
$pos = $p.Range.Start
$r = $d.Range($pos, $pos)
$r.InsertAfter('Referenser - knärot')
$pos = $pos + 19

# Paragraph 7: style=None
$p = $d.Paragraphs.Item($anchorIdx + 8)
$p.Style = $d.Styles.Item("Normal")
$pos = $p.Range.Start
$r = $d.Range($pos, $pos)
$r.InsertAfter('de Graaf M & Roberts M.R., 2009. ')
$pos = $pos + 33
$r = $d.Range($pos, $pos)
$r.InsertAfter('Short-term response of the herbaceous layer within leave patches after harvest. ')
$d.Range($pos, $pos + 80).Font.Italic = $true
$pos = $pos + 80
$r = $d.Range($pos, $pos)
$r.InsertAfter('Forest Ecology and Management 257, 1014-1025')
$pos = $pos + 44

# Paragraph 8: style=None
$p = $d.Paragraphs.Item($anchorIdx + 9)
$p.Style = $d.Styles.Item("Normal")
$pos = $p.Range.Start
$r = $d.Range($pos, $pos)
$r.InsertAfter('Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. ')
$pos = $pos + 62
$r = $d.Range($pos, $pos)
$r.InsertAfter('Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. ')
$d.Range($pos, $pos + 114).Font.Italic = $true
$pos = $pos + 114
$r = $d.Range($pos, $pos)
$r.InsertAfter('Ecological Applications, 22, 2049-2064 ')
$pos = $pos + 39

# Paragraph 9: style=None
$p = $d.Paragraphs.Item($anchorIdx + 10)
$p.Style = $d.Styles.Item("Normal")
$pos = $p.Range.Start
$r = $d.Range($pos, $pos)
$r.InsertAfter('Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. ')
$pos = $pos + 117
$r = $d.Range($pos, $pos)
$r.InsertAfter('Interactive effects of drought and edge exposure on old-growth forest understory species. ')
$d.Range($pos, $pos + 90).Font.Italic = $true
$pos = $pos + 90
$r = $d.Range($pos, $pos)
$r.InsertAfter('Landscape Ecology, 37, sid 1839-1853')
$pos = $pos + 36

# Paragraph 10: style=None
$p = $d.Paragraphs.Item($anchorIdx + 11)
$p.Style = $d.Styles.Item("Normal")
$pos = $p.Range.Start
$r = $d.Range($pos, $pos)
$r.InsertAfter('Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. ')
$pos = $pos + 54
$r = $d.Range($pos, $pos)
$r.InsertAfter('Biological legacies buffer local species extinction after logging. ')
$d.Range($pos, $pos + 67).Font.Italic = $true
$pos = $pos + 67
$r = $d.Range($pos, $pos)
$r.InsertAfter('Journal of Applied Ecology. 51, 53-62.')
$pos = $pos + 38

# Paragraph 11: style=None
$p = $d.Paragraphs.Item($anchorIdx + 12)
$p.Style = $d.Styles.Item("Normal")
$pos = $p.Range.Start
$r = $d.Range($pos, $pos)
$r.InsertAfter('Skogsstyrelsen, 2022. ')
$pos = $pos + 22
$r = $d.Range($pos, $pos)
$r.InsertAfter('Vägledning för hänsyn till knärot. ')
$d.Range($pos, $pos + 35).Font.Italic = $true
$pos = $pos + 35
$r = $d.Range($pos, $pos)
$r.InsertAfter('https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/')
$pos = $pos + 128

# Paragraph 12: style=None
$p = $d.Paragraphs.Item($anchorIdx + 13)
$p.Style = $d.Styles.Item("Normal")
$pos = $p.Range.Start
$r = $d.Range($pos, $pos)
$r.InsertAfter('SLU Artdatabanken, 2021. ')
$pos = $pos + 25
$r = $d.Range($pos, $pos)
$r.InsertAfter('Artfaktablad. Naturvård – artfakta. ')
$d.Range($pos, $pos + 36).Font.Italic = $true
$pos = $pos + 36
$r = $d.Range($pos, $pos)
$r.InsertAfter('SLU Artdatabanken, Uppsala ')
$pos = $pos + 27

# Update the report date shown in the first-page header.
$hdr = $d.Sections.Item(1).Headers.Item(2)
$hdr.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2) | Out-Null
